$wb = $excel.ActiveWorkbook

# Both "展览" and "全部类型" sheets carry the same exhibition rows and both
# need their "想去人数" (want-to-go count) figures bumped for rows 3 and 4.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 74
    $ws.Range("F4").Value = 46
}
